$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.596.27'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.718.87'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('D4').Value = '0.9979'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '240.54'
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('D6').Value = '0.9987'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('D8').Value = '0.2598'
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('D9').Value = '0.06204'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').Value = '1.726.61'
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('D11').Value = '0.07001'
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').Value = '15.74'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '0.6062'
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('D14').Value = '4.484'
$ws.Range('E14').Value = '  -2.50%  '
$ws.Range('D15').Value = '76.77'
$ws.Range('E15').Value = '  -1.64%  '
$ws.Range('D16').Value = '0.9982'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').Value = '26.434.68'
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').Value = '0.9981'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = '0.000007146'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').Value = '1.943.78'
$ws.Range('E21').Value = '  -1.18%  '
$ws.Range('D22').Value = '4.409'
$ws.Range('E22').Value = '  -3.32%  '
$ws.Range('D23').Value = '8.509'
$ws.Range('E23').Value = '  -2.57%  '
$ws.Range('D24').Value = '5.075'
$ws.Range('E24').Value = '  -4.08%  '
$ws.Range('D25').Value = '137.70'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').Value = '15.28'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').Value = '1.404'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').Value = '1.735'
$ws.Range('E28').Value = '  -1.42%  '
$ws.Range('D29').Value = '105.48'
$ws.Range('D30').Value = '3.915'
$ws.Range('E30').Value = '  -2.63%  '
$ws.Range('D31').Value = '0.07948'
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('D32').Value = '3.641'
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('D33').Value = '0.04498'
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('D34').Value = '2.608'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').Value = '0.9979'
$ws.Range('E35').Value = '  -1.83%  '
$ws.Range('D36').Value = '0.6228'
$ws.Range('E36').Value = '  -2.22%  '
$ws.Range('D37').Value = '0.9354'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('D38').Value = '1.996'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('D39').Value = '2.407'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').Value = '0.9985'
$ws.Range('D41').Value = '0.01511'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '5.515'
$ws.Range('E42').Value = '  +1.24%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '99.20'
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('D44').Value = '0.3832'
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('D45').Value = '6.904'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').Value = '0.1155'
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').Value = '0.05370'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').Value = '7.757'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('D49').Value = '30.07'
$ws.Range('E49').Value = '  -2.21%  '
$ws.Range('D50').Value = '51.37'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('E51').Value = '  -2.42%  '
